$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Supplier Part" header in column O (15)
$ws.Cells.Item(1, 15).Value = "Supplier Part"

# Fix the Comment (column N) for row 11: "pi" -> "Pi Header"
$ws.Cells.Item(11, 14).Value = "Pi Header"

# Fill in the new Supplier Part column (O) values for each data row, in order
$supplierParts = @(
    "C124378",   # row 2  - J2
    "C124378",   # row 3  - J2
    "C96446",    # row 4  - C7
    "C96446",    # row 5  - C8
    "C14663",    # row 6  - C9
    "C14663",    # row 7  - C10
    "C5120765",  # row 8  - U4
    "C28323",    # row 9  - C3
    "C16133",    # row 10 - C1
    "C2977589",  # row 11 - U1
    "C411294",   # row 12 - U3
    "C914555",   # row 13 - RF1
    "C11702",    # row 14 - R1
    "C11702",    # row 15 - R2
    "C86038",    # row 16 - C2
    "C86038",    # row 17 - C4
    "",          # row 18 - U5 (no supplier part)
    ""           # row 19 - J3 (no supplier part)
)

for ($i = 0; $i -lt $supplierParts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 15).Value = $supplierParts[$i]
}

# Make column O match the display width of the other columns (20 characters)
$ws.Columns.Item(15).ColumnWidth = 19.14
